# Update RAD Summary test case "Date" (column B) values across worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Mon Feb 24 23:07:42 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 23:08:01 EST 2025"
$ws.Range("B4").Value = "Mon Feb 24 23:08:20 EST 2025"
$ws.Range("B5").Value = "Mon Feb 24 23:08:39 EST 2025"
$ws.Range("B6").Value = "Mon Feb 24 23:08:58 EST 2025"
$ws.Range("B7").Value = "Mon Feb 24 23:09:17 EST 2025"

$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Mon Feb 24 23:09:39 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 23:09:59 EST 2025"
$ws.Range("B4").Value = "Mon Feb 24 23:10:19 EST 2025"
$ws.Range("B5").Value = "Mon Feb 24 23:10:39 EST 2025"
$ws.Range("B6").Value = "Mon Feb 24 23:10:57 EST 2025"
$ws.Range("B7").Value = "Mon Feb 24 23:11:16 EST 2025"
$ws.Range("B8").Value = "Mon Feb 24 23:11:39 EST 2025"
$ws.Range("B9").Value = "Mon Feb 24 23:11:59 EST 2025"
$ws.Range("B10").Value = "Mon Feb 24 23:12:20 EST 2025"
$ws.Range("B11").Value = "Mon Feb 24 23:12:41 EST 2025"
$ws.Range("B12").Value = "Mon Feb 24 23:13:00 EST 2025"
$ws.Range("B13").Value = "Mon Feb 24 23:13:21 EST 2025"
$ws.Range("B14").Value = "Mon Feb 24 23:13:43 EST 2025"
$ws.Range("B15").Value = "Mon Feb 24 23:14:03 EST 2025"
$ws.Range("B16").Value = "Mon Feb 24 23:14:22 EST 2025"
$ws.Range("B17").Value = "Mon Feb 24 23:14:43 EST 2025"
$ws.Range("B18").Value = "Mon Feb 24 23:15:02 EST 2025"
$ws.Range("B19").Value = "Mon Feb 24 23:15:20 EST 2025"
$ws.Range("B20").Value = "Mon Feb 24 23:15:41 EST 2025"
$ws.Range("B21").Value = "Mon Feb 24 23:16:00 EST 2025"

$ws = $wb.Worksheets.Item("Extension")
$ws.Range("B2").Value = "Mon Feb 24 23:16:20 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 23:16:42 EST 2025"
$ws.Range("B4").Value = "Mon Feb 24 23:17:00 EST 2025"
$ws.Range("B5").Value = "Mon Feb 24 23:17:19 EST 2025"
$ws.Range("B6").Value = "Mon Feb 24 23:17:38 EST 2025"
$ws.Range("B7").Value = "Mon Feb 24 23:17:58 EST 2025"

$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Mon Feb 24 23:18:17 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 23:18:37 EST 2025"
$ws.Range("B4").Value = "Mon Feb 24 23:18:57 EST 2025"
$ws.Range("B5").Value = "Mon Feb 24 23:19:16 EST 2025"
$ws.Range("B6").Value = "Mon Feb 24 23:19:35 EST 2025"
$ws.Range("B7").Value = "Mon Feb 24 23:19:55 EST 2025"
$ws.Range("B8").Value = "Mon Feb 24 23:20:15 EST 2025"
$ws.Range("B9").Value = "Mon Feb 24 23:20:36 EST 2025"
$ws.Range("B10").Value = "Mon Feb 24 23:20:57 EST 2025"
$ws.Range("B11").Value = "Mon Feb 24 23:21:17 EST 2025"
$ws.Range("B12").Value = "Mon Feb 24 23:21:37 EST 2025"
$ws.Range("B13").Value = "Mon Feb 24 23:21:57 EST 2025"
$ws.Range("B14").Value = "Mon Feb 24 23:22:16 EST 2025"
$ws.Range("B15").Value = "Mon Feb 24 23:22:37 EST 2025"
$ws.Range("B16").Value = "Mon Feb 24 23:22:57 EST 2025"
$ws.Range("B17").Value = "Mon Feb 24 23:23:18 EST 2025"
$ws.Range("B18").Value = "Mon Feb 24 23:23:42 EST 2025"
$ws.Range("B19").Value = "Mon Feb 24 23:24:07 EST 2025"
$ws.Range("B20").Value = "Mon Feb 24 23:24:32 EST 2025"
$ws.Range("B21").Value = "Mon Feb 24 23:24:54 EST 2025"
$ws.Range("B22").Value = "Mon Feb 24 23:25:15 EST 2025"
$ws.Range("B23").Value = "Mon Feb 24 23:25:38 EST 2025"
$ws.Range("B24").Value = "Mon Feb 24 23:26:02 EST 2025"
$ws.Range("B25").Value = "Mon Feb 24 23:26:25 EST 2025"
$ws.Range("B26").Value = "Mon Feb 24 23:26:47 EST 2025"
$ws.Range("B27").Value = "Mon Feb 24 23:27:09 EST 2025"
$ws.Range("B28").Value = "Mon Feb 24 23:27:31 EST 2025"
$ws.Range("B29").Value = "Mon Feb 24 23:27:54 EST 2025"
$ws.Range("B30").Value = "Mon Feb 24 23:28:15 EST 2025"
$ws.Range("B31").Value = "Mon Feb 24 23:28:38 EST 2025"
$ws.Range("B32").Value = "Mon Feb 24 23:28:59 EST 2025"
$ws.Range("B33").Value = "Mon Feb 24 23:29:19 EST 2025"
$ws.Range("B34").Value = "Mon Feb 24 23:29:41 EST 2025"
$ws.Range("B35").Value = "Mon Feb 24 23:30:02 EST 2025"
$ws.Range("B36").Value = "Mon Feb 24 23:30:20 EST 2025"
$ws.Range("B37").Value = "Mon Feb 24 23:30:42 EST 2025"
$ws.Range("B38").Value = "Mon Feb 24 23:31:03 EST 2025"
$ws.Range("B39").Value = "Mon Feb 24 23:31:25 EST 2025"
$ws.Range("B40").Value = "Mon Feb 24 23:31:44 EST 2025"
$ws.Range("B41").Value = "Mon Feb 24 23:32:06 EST 2025"
$ws.Range("B42").Value = "Mon Feb 24 23:32:27 EST 2025"
$ws.Range("B43").Value = "Mon Feb 24 23:32:48 EST 2025"
$ws.Range("B44").Value = "Mon Feb 24 23:33:08 EST 2025"
$ws.Range("B45").Value = "Mon Feb 24 23:33:28 EST 2025"
$ws.Range("B46").Value = "Mon Feb 24 23:33:50 EST 2025"
$ws.Range("B47").Value = "Mon Feb 24 23:34:12 EST 2025"
$ws.Range("B48").Value = "Mon Feb 24 23:34:32 EST 2025"
$ws.Range("B49").Value = "Mon Feb 24 23:34:54 EST 2025"
$ws.Range("B50").Value = "Mon Feb 24 23:35:15 EST 2025"
$ws.Range("B51").Value = "Mon Feb 24 23:35:36 EST 2025"
$ws.Range("B52").Value = "Mon Feb 24 23:35:59 EST 2025"
$ws.Range("B53").Value = "Mon Feb 24 23:36:20 EST 2025"
$ws.Range("B54").Value = "Mon Feb 24 23:36:40 EST 2025"
$ws.Range("B55").Value = "Mon Feb 24 23:37:00 EST 2025"
$ws.Range("B56").Value = "Mon Feb 24 23:37:20 EST 2025"
$ws.Range("B57").Value = "Mon Feb 24 23:37:41 EST 2025"
$ws.Range("B58").Value = "Mon Feb 24 23:38:02 EST 2025"
$ws.Range("B59").Value = "Mon Feb 24 23:38:25 EST 2025"
$ws.Range("B60").Value = "Mon Feb 24 23:38:47 EST 2025"
$ws.Range("B61").Value = "Mon Feb 24 23:39:07 EST 2025"

$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Mon Feb 24 23:40:14 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 23:40:36 EST 2025"
$ws.Range("B4").Value = "Mon Feb 24 23:40:59 EST 2025"
$ws.Range("B5").Value = "Mon Feb 24 23:41:22 EST 2025"
$ws.Range("B6").Value = "Mon Feb 24 23:41:46 EST 2025"
$ws.Range("B7").Value = "Mon Feb 24 23:42:09 EST 2025"
$ws.Range("B8").Value = "Mon Feb 24 23:42:35 EST 2025"
$ws.Range("B9").Value = "Mon Feb 24 23:42:58 EST 2025"

$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Range("B2").Value = "Mon Feb 24 23:43:27 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 23:44:02 EST 2025"
$ws.Range("B4").Value = "Mon Feb 24 23:44:35 EST 2025"
$ws.Range("B5").Value = "Mon Feb 24 23:45:22 EST 2025"
$ws.Range("B6").Value = "Mon Feb 24 23:45:57 EST 2025"

$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Mon Feb 24 23:39:29 EST 2025"
$ws.Range("B3").Value = "Mon Feb 24 23:39:49 EST 2025"
